$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "extra_column"
$ws.Range("G2").Value = "extra_column_data"
$ws.Range("G3").Select()
